# Apply the README/docx data-refresh edit described by the commit.
$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Simple single-value cell updates (1-indexed rows).
$tbl.Cell(1, 1).Range.Text  = "0M"
$tbl.Cell(2, 1).Range.Text  = "0M"
$tbl.Cell(3, 1).Range.Text  = "0M"
$tbl.Cell(4, 1).Range.Text  = "858"
$tbl.Cell(10, 1).Range.Text = "0.00006"
$tbl.Cell(11, 1).Range.Text = "0.00009"
$tbl.Cell(12, 1).Range.Text = "0.04606"

# Rows 44-46 previously contained a whole tab-separated data dump crammed
# into a single run; replace each with just the single correct value.
$tbl.Cell(44, 1).Range.Text = "100"
$tbl.Cell(45, 1).Range.Text = "0.05"
$tbl.Cell(46, 1).Range.Text = "3854"
